# Auto-generated edit script: apply scheduled-runner value updates
# across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the
# canonical-OOXML diff for this commit.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 559.6667
$ws.Range("I58").Value = 89.5
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 268.5
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -118.5
$ws.Range("N58").Value = -4800
$ws.Range("H86").Value = 5527.4165
$ws.Range("J86").Value = 6466.3
$ws.Range("L86").Value = 6466.3
$ws.Range("N86").Value = -8712.299999999999
$ws.Range("H89").Value = 5527.4165
$ws.Range("J89").Value = 6466.3
$ws.Range("L89").Value = 32331.5
$ws.Range("N89").Value = -43563.5
$ws.Range("H94").Value = 764.6667
$ws.Range("I94").Value = 764.6667
$ws.Range("K94").Value = 764.6667
$ws.Range("M94").Value = -313.6667
$ws.Range("H132").Value = 1027.9166
$ws.Range("J132").Value = 999.5
$ws.Range("L132").Value = 2998.5
$ws.Range("N132").Value = -8058.5
$ws.Range("H138").Value = 3509.4487
$ws.Range("J138").Value = 3450.0476
$ws.Range("L138").Value = 10350.1428
$ws.Range("N138").Value = -20630.1428

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4170.4443
$ws.Range("I32").Value = 3709.9119
$ws.Range("J32").Value = 11999.5
$ws.Range("K32").Value = 3709.9119
$ws.Range("L32").Value = 11999.5
$ws.Range("M32").Value = -3422.9119
$ws.Range("N32").Value = -12573.5
$ws.Range("H45").Value = 1684.3077
$ws.Range("I45").Value = 1627
$ws.Range("K45").Value = 1627
$ws.Range("M45").Value = -1250
$ws.Range("H61").Value = 3614.2917
$ws.Range("I61").Value = 1383.4
$ws.Range("J61").Value = 7332.4443
$ws.Range("K61").Value = 1383.4
$ws.Range("L61").Value = 7332.4443
$ws.Range("M61").Value = -1171.4
$ws.Range("N61").Value = -7756.4443
$ws.Range("H74").Value = 2631.1333
$ws.Range("I74").Value = 2512.0715
$ws.Range("J74").Value = 4298
$ws.Range("K74").Value = 2512.0715
$ws.Range("L74").Value = 4298
$ws.Range("M74").Value = -1638.0715
$ws.Range("N74").Value = -6046
$ws.Range("H77").Value = 2631.1333
$ws.Range("I77").Value = 2512.0715
$ws.Range("J77").Value = 4298
$ws.Range("K77").Value = 12560.3575
$ws.Range("L77").Value = 21490
$ws.Range("M77").Value = -8192.3575
$ws.Range("N77").Value = -30226
$ws.Range("H132").Value = 2057.4546
$ws.Range("I132").Value = 1916.6666
$ws.Range("K132").Value = 5749.9998
$ws.Range("M132").Value = -3219.9998
$ws.Range("H136").Value = 3614.2917
$ws.Range("I136").Value = 1383.4
$ws.Range("J136").Value = 7332.4443
$ws.Range("K136").Value = 4150.200000000001
$ws.Range("L136").Value = 21997.3329
$ws.Range("M136").Value = -1600.200000000001
$ws.Range("N136").Value = -27097.3329

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1459
$ws.Range("I134").Value = 1526.9
$ws.Range("K134").Value = 4580.700000000001
$ws.Range("M134").Value = -2045.700000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4637.3335
$ws.Range("I16").Value = 4637.3335
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4637.3335
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4350.3335
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 3491.4666
$ws.Range("I31").Value = 3227
$ws.Range("J31").Value = 3623.7
$ws.Range("K31").Value = 3227
$ws.Range("L31").Value = 3623.7
$ws.Range("M31").Value = -2932
$ws.Range("N31").Value = -4213.7
$ws.Range("H34").Value = 3491.4666
$ws.Range("I34").Value = 3227
$ws.Range("J34").Value = 3623.7
$ws.Range("K34").Value = 3227
$ws.Range("L34").Value = 3623.7
$ws.Range("M34").Value = -3025
$ws.Range("N34").Value = -4027.7
$ws.Range("H58").Value = 2989.75
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H105").Value = 3215.7917
$ws.Range("I105").Value = 2248.4167
$ws.Range("J105").Value = 4183.1665
$ws.Range("K105").Value = 2248.4167
$ws.Range("L105").Value = 4183.1665
$ws.Range("M105").Value = -501.4167000000002
$ws.Range("N105").Value = -7677.1665
$ws.Range("H107").Value = 580.6875
$ws.Range("J107").Value = 918.2
$ws.Range("L107").Value = 918.2
$ws.Range("N107").Value = -4758.2
$ws.Range("H113").Value = 4637.3335
$ws.Range("I113").Value = 4637.3335
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4637.3335
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2467.3335
$ws.Range("N113").ClearContents()
$ws.Range("H134").Value = 2153.4075
$ws.Range("I134").Value = 2112.8696
$ws.Range("K134").Value = 6338.6088
$ws.Range("M134").Value = -3803.6088
$ws.Range("H136").Value = 2989.75
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 514.7143
$ws.Range("I122").Value = 513
$ws.Range("J122").Value = 516
$ws.Range("K122").Value = 4617
$ws.Range("L122").Value = 4644
$ws.Range("M122").Value = -2167
$ws.Range("N122").Value = -9544

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 35000
$ws.Range("J121").Value = 35000
$ws.Range("L121").Value = 35000
$ws.Range("N121").Value = -38494

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 454.14285
$ws.Range("I16").Value = 454.14285
$ws.Range("K16").Value = 454.14285
$ws.Range("M16").Value = -284.14285
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -826
$ws.Range("H22").Value = 3160.5715
$ws.Range("J22").Value = 2425.1428
$ws.Range("L22").Value = 2425.1428
$ws.Range("N22").Value = -3015.1428
$ws.Range("H27").Value = 3160.5715
$ws.Range("J27").Value = 2425.1428
$ws.Range("L27").Value = 2425.1428
$ws.Range("N27").Value = -2639.1428
$ws.Range("H32").Value = 15000
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20634
$ws.Range("H46").Value = 2876.8096
$ws.Range("I46").Value = 1745.6666
$ws.Range("K46").Value = 1745.6666
$ws.Range("M46").Value = -1557.6666

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000822.3
$ws.Range("I81").Value = 983
$ws.Range("J81").Value = 3333780.8
$ws.Range("K81").Value = 1966
$ws.Range("L81").Value = 6667561.6
$ws.Range("M81").Value = -905
$ws.Range("N81").Value = -6669683.6
$ws.Range("H84").Value = 1000822.3
$ws.Range("I84").Value = 983
$ws.Range("J84").Value = 3333780.8
$ws.Range("K84").Value = 9830
$ws.Range("L84").Value = 33337808
$ws.Range("M84").Value = -4526
$ws.Range("N84").Value = -33348416
$ws.Range("H100").Value = 11112272
$ws.Range("I100").Value = 20000864
$ws.Range("K100").Value = 40001728
$ws.Range("M100").Value = -40001187
